# Regenerate the word / image / category cue table (Sheet1, rows 2-49)
# with the locked-in randomized assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cueData = @(
    @(2, "wenden", "none", "none"),
    @(3, "platzen", "dog/dog030.jpg", "dog"),
    @(4, "leuchten", "dog/dog009.jpg", "dog"),
    @(5, "rufen", "none", "none"),
    @(6, "fragen", "dog/dog000.jpg", "dog"),
    @(7, "schulden", "face/face021.jpg", "face"),
    @(8, "enden", "none", "none"),
    @(9, "tragen", "dog/dog022.jpg", "dog"),
    @(10, "streichen", "face/face015.jpg", "face"),
    @(11, "opfern", "none", "none"),
    @(12, "wohnen", "dog/dog012.jpg", "dog"),
    @(13, "rühren", "face/face010.jpg", "face"),
    @(14, "schalten", "none", "none"),
    @(15, "sammeln", "dog/dog010.jpg", "dog"),
    @(16, "faulen", "face/face005.jpg", "face"),
    @(17, "klagen", "none", "none"),
    @(18, "wehtun", "dog/dog029.jpg", "dog"),
    @(19, "gelten", "dog/dog007.jpg", "dog"),
    @(20, "sparen", "none", "none"),
    @(21, "schicken", "face/face023.jpg", "face"),
    @(22, "lassen", "dog/dog005.jpg", "dog"),
    @(23, "hören", "none", "none"),
    @(24, "treiben", "face/face025.jpg", "face"),
    @(25, "starren", "dog/dog021.jpg", "dog"),
    @(26, "drohen", "none", "none"),
    @(27, "herrschen", "dog/dog016.jpg", "dog"),
    @(28, "bilden", "dog/dog003.jpg", "dog"),
    @(29, "orten", "none", "none"),
    @(30, "runden", "dog/dog013.jpg", "dog"),
    @(31, "segeln", "dog/dog015.jpg", "dog"),
    @(32, "weigern", "none", "none"),
    @(33, "klingen", "face/face016.jpg", "face"),
    @(34, "schütteln", "face/face030.jpg", "face"),
    @(35, "ächzen", "none", "none"),
    @(36, "biegen", "face/face014.jpg", "face"),
    @(37, "öffnen", "dog/dog006.jpg", "dog"),
    @(38, "kosten", "none", "none"),
    @(39, "treten", "dog/dog014.jpg", "dog"),
    @(40, "küssen", "face/face017.jpg", "face"),
    @(41, "dauern", "none", "none"),
    @(42, "packen", "face/face003.jpg", "face"),
    @(43, "schaden", "face/face026.jpg", "face"),
    @(44, "stören", "none", "none"),
    @(45, "fügen", "face/face031.jpg", "face"),
    @(46, "wundern", "face/face008.jpg", "face"),
    @(47, "bremsen", "none", "none"),
    @(48, "reizen", "face/face027.jpg", "face"),
    @(49, "schultern", "face/face007.jpg", "face")
)

foreach ($cue in $cueData) {
    $r = $cue[0]
    $ws.Cells.Item($r, 1).Value = $cue[1]   # word
    $ws.Cells.Item($r, 2).Value = $cue[2]   # image
    $ws.Cells.Item($r, 3).Value = $cue[3]   # category
}

